# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the
# ad59fcdb-33cb-4bfb-a2ef-5492b3849380.md row, as produced by a fresh
# handback report generation run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to ad59fcdb-...md
# Column G = "Latest HO Xliff Generate Date"
$overview.Range("G3").Value = "2016-09-06 21:01:49"

# zh-cn sheet: row 3 corresponds to ad59fcdb-...md
# Column H = "Correspond Handoff Datetime", Column K = "Correspond Handback DateTime"
$zhcn.Range("H3").Value = "2016-09-06 21:01:44"
$zhcn.Range("K3").Value = "2016-09-06 21:02:34"

# de-de sheet: row 3 corresponds to ad59fcdb-...md
# Column K = "Correspond Handback DateTime"
$dede.Range("K3").Value = "2016-09-06 21:02:44"
